$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Himel (row 5) got a guest meal on day 10 (column K): 2.5 -> 5
$ws.Range("K5").Value = 5

# Bazar (row 43) Moricher gura (chili powder) purchase of 10 on day 10 (column K): 0 -> 10
$ws.Range("K43").Value = 10

# Update the active selection to K6
[void]$ws.Range("K6").Select()
